$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.55
$ws.Range("H2").Value = 4.5
$ws.Range("J2").Value = 2.05
$ws.Range("K2").Value = 2.5
$ws.Range("N2").Value = 17
$ws.Range("Q2").Value = 1.53
$ws.Range("R2").Value = 2.5
$ws.Range("S2").Value = 1.29
$ws.Range("T2").Value = 3.5
$ws.Range("U2").Value = 1.62
$ws.Range("V2").Value = 2.2
$ws.Range("W2").Value = 9.5
$ws.Range("X2").Value = 9
$ws.Range("AA2").Value = 11
$ws.Range("AC2").Value = 17
$ws.Range("AD2").Value = 9
$ws.Range("AH2").Value = 19
$ws.Range("AM2").Value = 34
$ws.Range("AO2").Value = 7.5
$ws.Range("AP2").Value = 15
$ws.Range("AQ2").Value = 21
$ws.Range("AS2").Value = 81
$ws.Range("AT2").Value = 3.5
$ws.Range("AU2").Value = 7.5
$ws.Range("AV2").Value = 41
$ws.Range("BA2").Value = 81
$ws.Range("BC2").Value = 401
$ws.Range("G3").Value = 2.3
$ws.Range("I3").Value = 3.4
$ws.Range("J3").Value = 3.1
$ws.Range("K3").Value = 1.95
$ws.Range("L3").Value = 4
$ws.Range("M3").Value = 1.11
$ws.Range("N3").Value = 6.5
$ws.Range("Q3").Value = 2.6
$ws.Range("R3").Value = 1.48
$ws.Range("U3").Value = 2.05
$ws.Range("V3").Value = 1.7
$ws.Range("X3").Value = 9.5
$ws.Range("AE3").Value = 17
$ws.Range("AI3").Value = 15
$ws.Range("AQ3").Value = 51
$ws.Range("AW3").Value = 5
$ws.Range("AZ3").Value = 67
$ws.Range("BA3").Value = 101
$ws.Range("BB3").Value = 301
$ws.Range("G4").Value = 3.3
$ws.Range("H4").Value = 3.1
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 3.1
$ws.Range("M4").Value = 1.08
$ws.Range("N4").Value = 8
$ws.Range("O4").Value = 1.4
$ws.Range("P4").Value = 2.75
$ws.Range("Q4").Value = 2.3
$ws.Range("R4").Value = 1.6
$ws.Range("S4").Value = 1.5
$ws.Range("T4").Value = 2.5
$ws.Range("U4").Value = 1.95
$ws.Range("V4").Value = 1.8
$ws.Range("W4").Value = 9
$ws.Range("AA4").Value = 29
$ws.Range("AB4").Value = 41
$ws.Range("AC4").Value = 8
$ws.Range("AG4").Value = 351
$ws.Range("AH4").Value = 7
$ws.Range("AI4").Value = 10
$ws.Range("AL4").Value = 21
$ws.Range("AM4").Value = 34
$ws.Range("AQ4").Value = 67
$ws.Range("AS4").Value = 251
$ws.Range("AT4").Value = 2.5
$ws.Range("AU4").Value = 8.5
$ws.Range("AV4").Value = 67
$ws.Range("AY4").Value = 26
$ws.Range("BB4").Value = 201
$ws.Range("N5").Value = 9
$ws.Range("Q5").Value = 2.25
$ws.Range("R5").Value = 1.62
